$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" --------------------------------------------
# Row 2 = CONSTANTE CAMACHO ARIANA ELIZABETH, col M = PORCELANATO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = -217.73

# --- Sheet "VENTA MENSUAL" ------------------------------------------------
# Row 2 = CONSTANTE CAMACHO ARIANA ELIZABETH, col F = junio
# Row 7 = totals row
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = -217.73
$ws2.Range("F7").Value = 633.6999999999999
# Column F widened (stored width 12 -> 13); ColumnWidth round-trips through
# character units, so compensate the ~0.8333 offset introduced on save.
$ws2.Columns.Item(6).ColumnWidth = 12.166666666666666

# --- Sheet "CUMPLIMIENTO MENSUAL" ----------------------------------------
# Row 14 = PORCELANATO, Row 17 = TOTAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D14").Value = -217.73
$ws3.Range("E14").Value = 7991.83
$ws3.Range("F14").Value = -0.02800710050037946

$ws3.Range("D17").Value = 633.6999999999999
$ws3.Range("E17").Value = 12866.3018254209
$ws3.Range("F17").Value = 0.04694073439358535

# Column D stored width 12 -> 13, column F stored width 18 -> 26
$ws3.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws3.Columns.Item(6).ColumnWidth = 25.166666666666668
